$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formats from column O into the new column P for the header/data rows,
# then set the actual values (this reuses existing style records, matching
# how the original author extended the table with a 2022 column).

$xlPasteFormats = -4122

# Row 4: new year header 2022
$ws.Range("O4").Copy() | Out-Null
$ws.Range("P4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("P4").Value = 2022

# Row 5
$ws.Range("O5").Copy() | Out-Null
$ws.Range("P5").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("P5").Value = 1

# Row 6
$ws.Range("O6").Copy() | Out-Null
$ws.Range("P6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("P6").Value = "-"

# Row 7
$ws.Range("O7").Copy() | Out-Null
$ws.Range("P7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("P7").Value = "-"

# Row 8
$ws.Range("O8").Copy() | Out-Null
$ws.Range("P8").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("P8").Value = "-"

# Row 9
$ws.Range("O9").Copy() | Out-Null
$ws.Range("P9").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("P9").Value = "-"

# Row 10
$ws.Range("O10").Copy() | Out-Null
$ws.Range("P10").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("P10").Value = "-"

# Row 11
$ws.Range("O11").Copy() | Out-Null
$ws.Range("P11").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("P11").Value = "-"

# Row 12
$ws.Range("O12").Copy() | Out-Null
$ws.Range("P12").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("P12").Value = 1

# Row 13
$ws.Range("O13").Copy() | Out-Null
$ws.Range("P13").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("P13").Value = "-"

# Row 14
$ws.Range("O14").Copy() | Out-Null
$ws.Range("P14").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("P14").Value = "-"

$excel.CutCopyMode = $false

# Move the active selection, matching the recorded end-state of the sheet.
$ws.Range("O21:O22").Select() | Out-Null
